# Generate Report for Handoff
# Updates Priority from "low" to "ht" and refreshes the "Latest Handoff Datetime"
# / "Latest HO Xliff Generate Date" timestamps for the rows that were ready for
# handoff, on the zh-cn and de-de worksheets (rows 4-7).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# zh-cn: rows 4-7 -> Priority (E) low -> ht, Latest Handoff Datetime (H) updated
for ($r = 4; $r -le 7; $r++) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-08-29 00:31:29"
}

# de-de: rows 4-7 -> Priority (E) low -> ht, Latest Handoff Datetime updated
for ($r = 4; $r -le 7; $r++) {
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-08-29 00:31:34"
}

# Overview: rows 4-7 -> Latest HO Xliff Generate Date (G) updated to match de-de
for ($r = 4; $r -le 7; $r++) {
    $overview.Cells.Item($r, 7).Value = "2016-08-29 00:31:34"
}
